$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $escaped = $val -replace '""', '""""'
    $ws.Range($cellRef).Formula = "=""" + $escaped + """"
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" "65.843.87"
Set-TextValue "E2" "  -3.33%  "
Set-TextValue "D3" "3.424.11"
Set-TextValue "E3" "  -5.27%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "184.27"
Set-TextValue "E5" "  -9.92%  "
Set-TextValue "D6" "534.83"
Set-TextValue "E6" "  -6.01%  "
Set-TextValue "D7" "0.617"
Set-TextValue "E7" "  -0.63%  "
Set-TextValue "D8" "3.413.75"
Set-TextValue "E8" "  -5.41%  "
Set-TextValue "E9" "  -0.05%  "
Set-TextValue "D10" "0.634"
Set-TextValue "E10" "  -6.15%  "
Set-TextValue "D11" "58.44"
Set-TextValue "E11" "  -4.12%  "
Set-TextValue "E12" "  -10.50%  "
Set-TextValue "D13" "0.0000259"
Set-TextValue "E13" "  -10.45%  "
Set-TextValue "D14" "9.47"
Set-TextValue "E14" "  -5.68%  "
Set-TextValue "D15" "3.959.16"
Set-TextValue "E15" "  -5.56%  "
Set-TextValue "E16" "  -2.50%  "
Set-TextValue "D17" "3.419.29"
Set-TextValue "E17" "  -5.43%  "
Set-TextValue "D18" "65.513.09"
Set-TextValue "E18" "  -3.55%  "
Set-TextValue "D19" "17.79"
Set-TextValue "E19" "  -6.01%  "
Set-TextValue "D20" "11.41"
Set-TextValue "E20" "  -7.95%  "
Set-TextValue "D21" "0.994"
Set-TextValue "E21" "  -7.82%  "
Set-TextValue "D22" "382.21"
Set-TextValue "E22" "  -5.27%  "
Set-TextValue "D23" "83.81"
Set-TextValue "E23" "  -2.00%  "
Set-TextValue "E24" "  -9.36%  "
Set-TextValue "D25" "10.94"
Set-TextValue "E25" "  -16.15%  "
Set-TextValue "B26" "InternetComputer(DFINITY)"
Set-TextValue "C26" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D26" "11.74"
Set-TextValue "E26" "  -7.25%  "
Set-TextValue "D27" "2.70"
Set-TextValue "E27" "  -8.65%  "
Set-TextValue "B28" "Toncoin"
Set-TextValue "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "3.66"
Set-TextValue "E28" "  -7.37%  "
Set-TextValue "D29" "8.64"
Set-TextValue "E29" "  -8.68%  "
Set-TextValue "D30" "689.03"
Set-TextValue "E30" "  +1.55%  "
Set-TextValue "D31" "30.16"
Set-TextValue "E31" "  -4.99%  "
Set-TextValue "D32" "6.80"
Set-TextValue "E32" "  -19.79%  "
Set-TextValue "B33" "Cosmos"
Set-TextValue "C33" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D33" "11.32"
Set-TextValue "E33" "  -7.58%  "
Set-TextValue "B34" "OKB"
Set-TextValue "C34" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D34" "61.92"
Set-TextValue "E34" "  -3.18%  "
Set-TextValue "E35" "  -6.44%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.30%  "
Set-TextValue "D37" "37.12"
Set-TextValue "E37" "  -12.24%  "
Set-TextValue "D38" "0.391"
Set-TextValue "E38" "  -7.69%  "
Set-TextValue "D39" "0.999"
Set-TextValue "E39" "  +0.13%  "
Set-TextValue "E40" "  -6.11%  "
Set-TextValue "D41" "2.901.97"
Set-TextValue "E41" "  -10.97%  "
Set-TextValue "D42" "2.79"
Set-TextValue "E42" "  -13.41%  "
Set-TextValue "D43" "2.71"
Set-TextValue "E43" "  -0.87%  "
Set-TextValue "B44" "PEPE"
Set-TextValue "C44" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D44" "0.0₃0634"
Set-TextValue "E44" "  -18.03%  "
Set-TextValue "B45" "VeChain"
Set-TextValue "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0396"
Set-TextValue "E45" "  -5.75%  "
Set-TextValue "D46" "2.36"
Set-TextValue "E46" "  -15.41%  "
Set-TextValue "E47" "  -3.52%  "
Set-TextValue "D48" "134.74"
Set-TextValue "E48" "  -3.76%  "
Set-TextValue "D49" "2.86"
Set-TextValue "E49" "  -6.74%  "
Set-TextValue "D50" "2.59"
Set-TextValue "E50" "  -5.37%  "
Set-TextValue "D51" "2.35"

$excel.CutCopyMode = 0

